# Update real estate listings: correct the bathroom count for row 7,
# replace the stale duplicate listing in row 9 (and re-align the rows that
# follow it), and fix the rounded price in row 12 - net effect matches the
# newly scraped data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 1.5

$ws.Range("A9").Value = "R 2 800 000"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "154 m²"
$ws.Range("E9").Value = "Kenilworth"

$ws.Range("A12").Value = "R 3 300 000"

$ws.Range("A13").Value = "R 3 450 000"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "95 m²"

$ws.Range("A14").Value = "R 3 595 000"
$ws.Range("D14").Value = "147 m²"

$ws.Range("A15").Value = "R 3 695 000"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = "200 m²"
$ws.Range("E15").Value = "Kenilworth Upper"

$ws.Range("A16").Value = "R 4 300 000"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 3.5
$ws.Range("D16").Value = "260 m²"
$ws.Range("E16").Value = "Lakeside"

$ws.Range("A17").Value = "R 4 800 000"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = "115 m²"
$ws.Range("E17").Value = "Newlands"

$ws.Range("A18").Value = "R 5 295 000"
$ws.Range("C18").Value = 3.5
$ws.Range("D18").Value = "314 m²"
$ws.Range("E18").Value = "Claremont"

$ws.Range("A19").Value = "R 8 300 000"
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = "154 m²"
$ws.Range("E19").Value = "Green Point"

$ws.Range("A20").Value = "R 999 999"
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 1.5
$ws.Range("D20").Value = "72 m²"
$ws.Range("E20").Value = "Thornton"

$ws.Range("A21").Value = "R 999 999"

$ws.Range("A22").Value = "R 999 999"
